$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '43.023.76'
$ws.Range('E2').Value = '  +0.96%  '
$ws.Range('D3').Value = '2.579.70'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.575'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.538'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.30'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.59'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('D13').Value = '2.975.78'
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.53%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.586.90'
$ws.Range('E16').Value = '  +2.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.847'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').Value = '43.079.13'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '0.0₃0971'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.19'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.94'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.45'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.81%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  +3.71%  '
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.93'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.54'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.97%  '
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.39%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.119'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.06'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.46%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.28'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.15%  '
$ws.Range('D46').Value = '2.010.54'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D48').Value = '2.825.57'
$ws.Range('E48').Value = '  +2.32%  '
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '82.08'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.54%  '
